$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1:F47")
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("C2:C47"), 0, 1, $null, 0)
$sort.SortFields.Add($ws.Range("E2:E47"), 0, 1, $null, 0)
$sort.SetRange($rng)
$sort.Header = 1
$sort.Apply()
Write-Output "sorted-via-sortobject"
